$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A473:A485").NumberFormat = "@"
$ws.Cells.Item(473, 1).Value = "2026-02-06"
$ws.Cells.Item(473, 2).Value = "10:19:15"
$ws.Cells.Item(473, 3).Value = "10:00"
$ws.Cells.Item(473, 4).Value = "Bathroom"
$ws.Cells.Item(473, 5).Value = "No Motion"
$ws.Cells.Item(473, 6).Value = "Inactive"
$ws.Cells.Item(474, 1).Value = "2026-02-06"
$ws.Cells.Item(474, 2).Value = "10:19:18"
$ws.Cells.Item(474, 3).Value = "10:00"
$ws.Cells.Item(474, 4).Value = "Bathroom"
$ws.Cells.Item(474, 5).Value = "No Motion"
$ws.Cells.Item(474, 6).Value = "Inactive"
$ws.Cells.Item(475, 1).Value = "2026-02-06"
$ws.Cells.Item(475, 2).Value = "10:19:21"
$ws.Cells.Item(475, 3).Value = "10:00"
$ws.Cells.Item(475, 4).Value = "Bathroom"
$ws.Cells.Item(475, 5).Value = "No Motion"
$ws.Cells.Item(475, 6).Value = "Inactive"
$ws.Cells.Item(476, 1).Value = "2026-02-06"
$ws.Cells.Item(476, 2).Value = "10:19:24"
$ws.Cells.Item(476, 3).Value = "10:00"
$ws.Cells.Item(476, 4).Value = "Bathroom"
$ws.Cells.Item(476, 5).Value = "No Motion"
$ws.Cells.Item(476, 6).Value = "Inactive"
$ws.Cells.Item(477, 1).Value = "2026-02-06"
$ws.Cells.Item(477, 2).Value = "10:19:30"
$ws.Cells.Item(477, 3).Value = "10:00"
$ws.Cells.Item(477, 4).Value = "Bathroom"
$ws.Cells.Item(477, 5).Value = "No Motion"
$ws.Cells.Item(477, 6).Value = "Inactive"
$ws.Cells.Item(478, 1).Value = "2026-02-06"
$ws.Cells.Item(478, 2).Value = "10:19:35"
$ws.Cells.Item(478, 3).Value = "10:00"
$ws.Cells.Item(478, 4).Value = "Bathroom"
$ws.Cells.Item(478, 5).Value = "No Motion"
$ws.Cells.Item(478, 6).Value = "Inactive"
$ws.Cells.Item(479, 1).Value = "2026-02-06"
$ws.Cells.Item(479, 2).Value = "10:19:40"
$ws.Cells.Item(479, 3).Value = "10:00"
$ws.Cells.Item(479, 4).Value = "Bathroom"
$ws.Cells.Item(479, 5).Value = "No Motion"
$ws.Cells.Item(479, 6).Value = "Inactive"
$ws.Cells.Item(480, 1).Value = "2026-02-06"
$ws.Cells.Item(480, 2).Value = "10:19:45"
$ws.Cells.Item(480, 3).Value = "10:00"
$ws.Cells.Item(480, 4).Value = "Bathroom"
$ws.Cells.Item(480, 5).Value = "No Motion"
$ws.Cells.Item(480, 6).Value = "Inactive"
$ws.Cells.Item(481, 1).Value = "2026-02-06"
$ws.Cells.Item(481, 2).Value = "10:19:50"
$ws.Cells.Item(481, 3).Value = "10:00"
$ws.Cells.Item(481, 4).Value = "Bathroom"
$ws.Cells.Item(481, 5).Value = "No Motion"
$ws.Cells.Item(481, 6).Value = "Inactive"
$ws.Cells.Item(482, 1).Value = "2026-02-06"
$ws.Cells.Item(482, 2).Value = "10:19:55"
$ws.Cells.Item(482, 3).Value = "10:00"
$ws.Cells.Item(482, 4).Value = "Bathroom"
$ws.Cells.Item(482, 5).Value = "No Motion"
$ws.Cells.Item(482, 6).Value = "Inactive"
$ws.Cells.Item(483, 1).Value = "2026-02-06"
$ws.Cells.Item(483, 2).Value = "10:20:00"
$ws.Cells.Item(483, 3).Value = "10:00"
$ws.Cells.Item(483, 4).Value = "Bathroom"
$ws.Cells.Item(483, 5).Value = "No Motion"
$ws.Cells.Item(483, 6).Value = "Inactive"
$ws.Cells.Item(484, 1).Value = "2026-02-06"
$ws.Cells.Item(484, 2).Value = "10:20:05"
$ws.Cells.Item(484, 3).Value = "10:00"
$ws.Cells.Item(484, 4).Value = "Bathroom"
$ws.Cells.Item(484, 5).Value = "No Motion"
$ws.Cells.Item(484, 6).Value = "Inactive"
$ws.Cells.Item(485, 1).Value = "2026-02-06"
$ws.Cells.Item(485, 2).Value = "10:20:10"
$ws.Cells.Item(485, 3).Value = "10:00"
$ws.Cells.Item(485, 4).Value = "Bathroom"
$ws.Cells.Item(485, 5).Value = "No Motion"
$ws.Cells.Item(485, 6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A327:A336").NumberFormat = "@"
$ws.Range("E327:E336").NumberFormat = "@"
$ws.Cells.Item(327, 1).Value = "2026-02-06"
$ws.Cells.Item(327, 2).Value = "10:19:16"
$ws.Cells.Item(327, 3).Value = "10:00"
$ws.Cells.Item(327, 4).Value = "Bathroom"
$ws.Cells.Item(327, 5).Value = "68.5%"
$ws.Cells.Item(327, 6).Value = "Active"
$ws.Cells.Item(328, 1).Value = "2026-02-06"
$ws.Cells.Item(328, 2).Value = "10:19:19"
$ws.Cells.Item(328, 3).Value = "10:00"
$ws.Cells.Item(328, 4).Value = "Bathroom"
$ws.Cells.Item(328, 5).Value = "68.5%"
$ws.Cells.Item(328, 6).Value = "Active"
$ws.Cells.Item(329, 1).Value = "2026-02-06"
$ws.Cells.Item(329, 2).Value = "10:19:22"
$ws.Cells.Item(329, 3).Value = "10:00"
$ws.Cells.Item(329, 4).Value = "Bathroom"
$ws.Cells.Item(329, 5).Value = "68.5%"
$ws.Cells.Item(329, 6).Value = "Active"
$ws.Cells.Item(330, 1).Value = "2026-02-06"
$ws.Cells.Item(330, 2).Value = "10:19:26"
$ws.Cells.Item(330, 3).Value = "10:00"
$ws.Cells.Item(330, 4).Value = "Bathroom"
$ws.Cells.Item(330, 5).Value = "68.4%"
$ws.Cells.Item(330, 6).Value = "Active"
$ws.Cells.Item(331, 1).Value = "2026-02-06"
$ws.Cells.Item(331, 2).Value = "10:19:31"
$ws.Cells.Item(331, 3).Value = "10:00"
$ws.Cells.Item(331, 4).Value = "Bathroom"
$ws.Cells.Item(331, 5).Value = "68.4%"
$ws.Cells.Item(331, 6).Value = "Active"
$ws.Cells.Item(332, 1).Value = "2026-02-06"
$ws.Cells.Item(332, 2).Value = "10:19:36"
$ws.Cells.Item(332, 3).Value = "10:00"
$ws.Cells.Item(332, 4).Value = "Bathroom"
$ws.Cells.Item(332, 5).Value = "68.4%"
$ws.Cells.Item(332, 6).Value = "Active"
$ws.Cells.Item(333, 1).Value = "2026-02-06"
$ws.Cells.Item(333, 2).Value = "10:19:41"
$ws.Cells.Item(333, 3).Value = "10:00"
$ws.Cells.Item(333, 4).Value = "Bathroom"
$ws.Cells.Item(333, 5).Value = "68.4%"
$ws.Cells.Item(333, 6).Value = "Active"
$ws.Cells.Item(334, 1).Value = "2026-02-06"
$ws.Cells.Item(334, 2).Value = "10:20:01"
$ws.Cells.Item(334, 3).Value = "10:00"
$ws.Cells.Item(334, 4).Value = "Bathroom"
$ws.Cells.Item(334, 5).Value = "68.0%"
$ws.Cells.Item(334, 6).Value = "Active"
$ws.Cells.Item(335, 1).Value = "2026-02-06"
$ws.Cells.Item(335, 2).Value = "10:20:06"
$ws.Cells.Item(335, 3).Value = "10:00"
$ws.Cells.Item(335, 4).Value = "Bathroom"
$ws.Cells.Item(335, 5).Value = "67.1%"
$ws.Cells.Item(335, 6).Value = "Active"
$ws.Cells.Item(336, 1).Value = "2026-02-06"
$ws.Cells.Item(336, 2).Value = "10:20:11"
$ws.Cells.Item(336, 3).Value = "10:00"
$ws.Cells.Item(336, 4).Value = "Bathroom"
$ws.Cells.Item(336, 5).Value = "68.1%"
$ws.Cells.Item(336, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A327:A336").NumberFormat = "@"
$ws.Cells.Item(327, 1).Value = "2026-02-06"
$ws.Cells.Item(327, 2).Value = "10:19:17"
$ws.Cells.Item(327, 3).Value = "10:00"
$ws.Cells.Item(327, 4).Value = "Bathroom"
$ws.Cells.Item(327, 5).Value = "28.2C"
$ws.Cells.Item(327, 6).Value = "Active"
$ws.Cells.Item(328, 1).Value = "2026-02-06"
$ws.Cells.Item(328, 2).Value = "10:19:20"
$ws.Cells.Item(328, 3).Value = "10:00"
$ws.Cells.Item(328, 4).Value = "Bathroom"
$ws.Cells.Item(328, 5).Value = "28.2C"
$ws.Cells.Item(328, 6).Value = "Active"
$ws.Cells.Item(329, 1).Value = "2026-02-06"
$ws.Cells.Item(329, 2).Value = "10:19:23"
$ws.Cells.Item(329, 3).Value = "10:00"
$ws.Cells.Item(329, 4).Value = "Bathroom"
$ws.Cells.Item(329, 5).Value = "28.3C"
$ws.Cells.Item(329, 6).Value = "Active"
$ws.Cells.Item(330, 1).Value = "2026-02-06"
$ws.Cells.Item(330, 2).Value = "10:19:27"
$ws.Cells.Item(330, 3).Value = "10:00"
$ws.Cells.Item(330, 4).Value = "Bathroom"
$ws.Cells.Item(330, 5).Value = "28.2C"
$ws.Cells.Item(330, 6).Value = "Active"
$ws.Cells.Item(331, 1).Value = "2026-02-06"
$ws.Cells.Item(331, 2).Value = "10:19:32"
$ws.Cells.Item(331, 3).Value = "10:00"
$ws.Cells.Item(331, 4).Value = "Bathroom"
$ws.Cells.Item(331, 5).Value = "28.2C"
$ws.Cells.Item(331, 6).Value = "Active"
$ws.Cells.Item(332, 1).Value = "2026-02-06"
$ws.Cells.Item(332, 2).Value = "10:19:37"
$ws.Cells.Item(332, 3).Value = "10:00"
$ws.Cells.Item(332, 4).Value = "Bathroom"
$ws.Cells.Item(332, 5).Value = "28.2C"
$ws.Cells.Item(332, 6).Value = "Active"
$ws.Cells.Item(333, 1).Value = "2026-02-06"
$ws.Cells.Item(333, 2).Value = "10:19:42"
$ws.Cells.Item(333, 3).Value = "10:00"
$ws.Cells.Item(333, 4).Value = "Bathroom"
$ws.Cells.Item(333, 5).Value = "28.3C"
$ws.Cells.Item(333, 6).Value = "Active"
$ws.Cells.Item(334, 1).Value = "2026-02-06"
$ws.Cells.Item(334, 2).Value = "10:20:02"
$ws.Cells.Item(334, 3).Value = "10:00"
$ws.Cells.Item(334, 4).Value = "Bathroom"
$ws.Cells.Item(334, 5).Value = "28.3C"
$ws.Cells.Item(334, 6).Value = "Active"
$ws.Cells.Item(335, 1).Value = "2026-02-06"
$ws.Cells.Item(335, 2).Value = "10:20:07"
$ws.Cells.Item(335, 3).Value = "10:00"
$ws.Cells.Item(335, 4).Value = "Bathroom"
$ws.Cells.Item(335, 5).Value = "28.3C"
$ws.Cells.Item(335, 6).Value = "Active"
$ws.Cells.Item(336, 1).Value = "2026-02-06"
$ws.Cells.Item(336, 2).Value = "10:20:12"
$ws.Cells.Item(336, 3).Value = "10:00"
$ws.Cells.Item(336, 4).Value = "Bathroom"
$ws.Cells.Item(336, 5).Value = "28.3C"
$ws.Cells.Item(336, 6).Value = "Active"
